# PHOENIX-5928: Completed the creation of journal voucher with budget check
#
# Adds a new "budgetCheck" journal-voucher test-data row to the
# journalVoucherDetails sheet, then makes that sheet/cell the active
# selection (it becomes the workbook's active tab).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("journalVoucherDetails")

# New row 5: budgetCheck voucher data
$ws.Range("A5").Value = "budgetCheck"
$ws.Range("B5").Value = "03/01/2016"
$ws.Range("C5").Value = "Expense"
$ws.Range("D5").Value = "2101001"
$ws.Range("E5").Value = "3501003"
$ws.Range("F5").Value = "ENGINEERING"
$ws.Range("G5").Value = "Water Supply"

# journalVoucherDetails becomes the active sheet, with G5 selected
$ws.Activate()
$ws.Range("G5").Select()
